# extraSpaceInEndBookmark: the field-result run group gets a fresh rsidR token
# (mirrors the regenerated revision-save id emitted by the updated parser /
# TokenIteratorFieldRewriterSplit) and the "missing ENDBOOKMARK" error text now
# gets a leading "    <---" marker so the problem is easier to spot in the
# generated document.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 2: "Test link before bookmark : <REF field>" -----------------
# Rebuild the paragraph verbatim, only swapping the w:rsidR token stamped on
# each run of the REF field (begin / instrText / separate / result / end).
$para2 = $d.Paragraphs(2)
$xmlPara2 = '<w:p ' + $wNs + ' w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidRPr="00FF681D">' +
  '<w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00FF681D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Test link before bookmark&#160;: </w:t></w:r>' +
  '<w:r w:rsidR="4DEB65CDE1DC46B2AD5084BE6FF020EB"><w:fldChar w:fldCharType="begin"/></w:r>' +
  '<w:r w:rsidR="4DEB65CDE1DC46B2AD5084BE6FF020EB"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r>' +
  '<w:r w:rsidR="4DEB65CDE1DC46B2AD5084BE6FF020EB"><w:fldChar w:fldCharType="separate"/></w:r>' +
  '<w:r w:rsidR="4DEB65CDE1DC46B2AD5084BE6FF020EB"><w:rPr><w:noProof/><w:b w:val="on"/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +
  '<w:r w:rsidR="4DEB65CDE1DC46B2AD5084BE6FF020EB"><w:fldChar w:fldCharType="end"/></w:r>' +
  '</w:p>'
$null = $para2.Range.InsertXML($xmlPara2)

# --- Paragraph 3: "Test bookmark : <bookmarked error text>" -----------------
# Rebuild the paragraph verbatim, bumping the bookmark id and prefixing the
# error message with the new "    <---" marker (now needs xml:space=preserve
# because of the leading spaces).
$para3 = $d.Paragraphs(3)
$quote = [char]39
$errText = '    &lt;---Invalid block: Unexpected tag EOF missing [ENDBOOKMARK] while parsing m:bookmark ' + $quote + 'bookmark1' + $quote
$xmlPara3 = '<w:p ' + $wNs + ' w:rsidP="009168BC" w:rsidR="00C52979" w:rsidRDefault="00E02A2B" w:rsidRPr="00FF681D">' +
  '<w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00FF681D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Test</w:t></w:r>' +
  '<w:r w:rsidR="00C52979" w:rsidRPr="00FF681D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r w:rsidRPr="00FF681D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bookmark</w:t></w:r>' +
  '<w:r w:rsidR="00C52979" w:rsidRPr="00FF681D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">&#160;: </w:t></w:r>' +
  '<w:bookmarkStart w:name="bookmark1" w:id="132473400371272052279446125497739496322"/>' +
  '<w:r><w:rPr><w:b w:val="on"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">' + $errText + '</w:t></w:r>' +
  '<w:bookmarkEnd w:id="132473400371272052279446125497739496322"/>' +
  '</w:p>'
$null = $para3.Range.InsertXML($xmlPara3)

Write-Host "applied"
